$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.990.97"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.304.70"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'252.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'0.641"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").Value = "'75.26"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.13%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.650"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").Value = "'38.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").Value = "'0.0991"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("D12").Value = "'7.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  +2.37%  "
$ws.Range("D14").Value = "2.652.66"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Value = "'15.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.38%  "
$ws.Range("D16").Value = "'0.881"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "2.316.78"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "42.884.54"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Value = "'6.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'72.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").Value = "'237.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "'2.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.63%  "
$ws.Range("D24").Value = "'3.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'11.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'2.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "'2.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'167.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'21.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0852"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.03%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'6.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "'31.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").Value = "'4.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.27%  "
$ws.Range("D37").Value = "'4.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("D38").Value = "'0.0308"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("D39").Value = "'13.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.04%  "
$ws.Range("D40").Value = "'2.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").Value = "'5.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").Value = "'0.215"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.80%  "
$ws.Range("D43").Value = "'9.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").Value = "'61.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").Value = "'4.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").Value = "'105.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.44%  "
$ws.Range("D47").Value = "'0.102"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'4.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
